$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.762.90"
Set-TextValue "E2" "  +0.41%  "
Set-TextValue "D3" "1.649.88"
Set-TextValue "E3" "  +1.07%  "
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "214.84"
Set-TextValue "E5" "  +0.82%  "
Set-TextValue "E6" "  +2.15%  "
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "0.252"
Set-TextValue "E8" "  -0.26%  "
Set-TextValue "E9" "  +0.50%  "
Set-TextValue "D10" "19.22"
Set-TextValue "E10" "  +1.09%  "
Set-TextValue "E11" "  -0.13%  "
Set-TextValue "D12" "1.871.42"
Set-TextValue "E12" "  +0.56%  "
Set-TextValue "D13" "1.641.26"
Set-TextValue "E13" "  +0.41%  "
Set-TextValue "E14" "  +1.94%  "
Set-TextValue "D15" "0.534"
Set-TextValue "E15" "  +1.29%  "
Set-TextValue "D16" "65.89"
Set-TextValue "E16" "  +4.34%  "
Set-TextValue "D17" "26.761.45"
Set-TextValue "E17" "  +0.38%  "
Set-TextValue "D18" "0.0₃0748"
Set-TextValue "E18" "  +1.05%  "
Set-TextValue "D19" "219.50"
Set-TextValue "E19" "  +4.43%  "
Set-TextValue "E20" "  -0.10%  "
Set-TextValue "D21" "4.36"
Set-TextValue "E21" "  +1.27%  "
Set-TextValue "D22" "6.34"
Set-TextValue "E22" "  +2.08%  "
Set-TextValue "D23" "9.45"
Set-TextValue "E23" "  -0.04%  "
Set-TextValue "E24" "  +11.64%  "
Set-TextValue "E25" "  +0.18%  "
Set-TextValue "E26" "  -0.09%  "
Set-TextValue "E27" "  -0.10%  "
Set-TextValue "D28" "6.96"
Set-TextValue "E28" "  +1.04%  "
Set-TextValue "D29" "15.84"
Set-TextValue "E29" "  +2.84%  "
Set-TextValue "E30" "  -0.41%  "
Set-TextValue "E31" "  -0.33%  "
Set-TextValue "E32" "  +4.20%  "
Set-TextValue "D33" "3.03"
Set-TextValue "E33" "  +2.68%  "
Set-TextValue "D34" "1.270.68"
Set-TextValue "E34" "  +8.62%  "
Set-TextValue "D35" "1.53"
Set-TextValue "E35" "  +1.28%  "
Set-TextValue "D36" "2.39"
Set-TextValue "E36" "  +1.07%  "
Set-TextValue "E37" "  +2.99%  "
Set-TextValue "D38" "0.814"
Set-TextValue "E38" "  +0.36%  "
Set-TextValue "D39" "0.515"
Set-TextValue "E39" "  +1.85%  "
Set-TextValue "E40" "  -0.13%  "
Set-TextValue "E41" "  -1.52%  "
Set-TextValue "E42" "  +1.02%  "
Set-TextValue "E43" "  -0.14%  "
Set-TextValue "D44" "1.782.71"
Set-TextValue "E44" "  +0.74%  "
Set-TextValue "D45" "93.95"
Set-TextValue "E45" "  +1.41%  "
Set-TextValue "E46" "  +3.80%  "
Set-TextValue "D47" "55.68"
Set-TextValue "E47" "  +1.92%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0101"
Set-TextValue "E48" "  -4.14%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.0514"
Set-TextValue "E49" "  +0.39%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "7.67"
Set-TextValue "E50" "  +1.15%  "
Set-TextValue "D51" "0.0968"
Set-TextValue "E51" "  +2.89%  "
